$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Window position/size (best effort - mirrors the author's window layout at save time)
$excel.ActiveWindow.Left = 1200
$excel.ActiveWindow.Top = 240
$excel.ActiveWindow.Width = 11025
$excel.ActiveWindow.Height = 10320

$ws.Range("B2").Value = "최준아"
$ws.Rows.Item(5).RowHeight = 33

$ws.Range("A21").Value = "rowbomb, sixbomb 제일 밑에서 `n쏘면 폭탄 소리 안나고 구슬 쌓이는 소리 남"
$ws.Range("A21").WrapText = $true
$ws.Range("B21").Value = "O"
$ws.Rows.Item(21).RowHeight = 49.5

$ws.Range("D3").Select()
